$d = $word.ActiveDocument

# 1. Replace the "Description" cell text for the weekly-plan requirement
#    with the new requirement wording (Use Case table, row "Lap ke hoach").
$oldText = "Chức năng này sẽ tạo một kế hoạch do chính User đặt ra mỗi tuần. User sẽ đặt ra cho mình số buổi tập luyện tối thiểu trong một tuần, hệ thống sẽ ghi nhận các kế hoạch mà User đặt ra và thông báo vào mỗi cuối tuần cho User có thể tiện giám sát quá trình tập luyện của mình."
$newText = "Chức năng cho phép user lập kế hoạch tập luyện hằng ngày. Giao diện có 1 lịch biểu theo tháng/ tuần, người dùng nhấp chọn ngày và chọn bộ phận muốn tập luyện để cải thiện"

$range = $d.Content
$range.Find.ClearFormatting()
[void]$range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# 2. The "Response" cell still carries a stray "_GoBack" bookmark (Word
#    re-marks the last edit position every save). Re-apply the same text
#    over itself so the save no longer anchors a bookmark there.
$range2 = $d.Content
$range2.Find.ClearFormatting()
[void]$range2.Find.Execute("Response", $false, $false, $false, $false, $false, $true, 1, $false, "Response", 2)
